# Update the MIBANCO mortgage rate table (table_mibanco.xlsx):
# within each "year" block of rows, the order of the (banks / rate_value)
# rows has been reversed, so the A (banks code) and E (rate_value)
# columns need to be updated to reflect the reversed order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 1).Value = 140
$ws.Cells.Item(3, 5).Value = 15
$ws.Cells.Item(4, 1).Value = 143
$ws.Cells.Item(4, 5).Value = 14.91
$ws.Cells.Item(5, 1).Value = 163
$ws.Cells.Item(5, 5).Value = 14.38
$ws.Cells.Item(6, 1).Value = 170
$ws.Cells.Item(6, 5).Value = 14.18
$ws.Cells.Item(7, 1).Value = 167
$ws.Cells.Item(7, 5).Value = 14.28
$ws.Cells.Item(8, 1).Value = 171
$ws.Cells.Item(8, 5).Value = 14.12
$ws.Cells.Item(10, 1).Value = 173
$ws.Cells.Item(10, 5).Value = 13.9
$ws.Cells.Item(11, 1).Value = 174
$ws.Cells.Item(11, 5).Value = 13.84
$ws.Cells.Item(12, 1).Value = 144
$ws.Cells.Item(12, 5).Value = 14.87
$ws.Cells.Item(13, 1).Value = 151
$ws.Cells.Item(13, 5).Value = 14.71
$ws.Cells.Item(15, 1).Value = 159
$ws.Cells.Item(15, 5).Value = 14.57
$ws.Cells.Item(16, 1).Value = 162
$ws.Cells.Item(16, 5).Value = 14.41
$ws.Cells.Item(17, 1).Value = 138
$ws.Cells.Item(17, 5).Value = 15.02
$ws.Cells.Item(18, 1).Value = 145
$ws.Cells.Item(18, 5).Value = 14.84
$ws.Cells.Item(19, 1).Value = 146
$ws.Cells.Item(20, 1).Value = 148
$ws.Cells.Item(20, 5).Value = 14.79
$ws.Cells.Item(21, 1).Value = 158
$ws.Cells.Item(21, 5).Value = 14.57
$ws.Cells.Item(22, 1).Value = 136
$ws.Cells.Item(22, 5).Value = 15.14
$ws.Cells.Item(23, 1).Value = 139
$ws.Cells.Item(23, 5).Value = 15
$ws.Cells.Item(25, 1).Value = 142
$ws.Cells.Item(25, 5).Value = 14.92
$ws.Cells.Item(26, 1).Value = 150
$ws.Cells.Item(26, 5).Value = 14.77
$ws.Cells.Item(27, 1).Value = 123
$ws.Cells.Item(27, 5).Value = 15.62
$ws.Cells.Item(28, 1).Value = 129
$ws.Cells.Item(28, 5).Value = 15.46
$ws.Cells.Item(30, 1).Value = 133
$ws.Cells.Item(30, 5).Value = 15.26
$ws.Cells.Item(31, 1).Value = 134
$ws.Cells.Item(31, 5).Value = 15.2
$ws.Cells.Item(32, 1).Value = 110
$ws.Cells.Item(32, 5).Value = 16.27
$ws.Cells.Item(33, 1).Value = 115
$ws.Cells.Item(33, 5).Value = 15.86
$ws.Cells.Item(35, 1).Value = 119
$ws.Cells.Item(35, 5).Value = 15.68
$ws.Cells.Item(36, 1).Value = 121
$ws.Cells.Item(36, 5).Value = 15.66
$ws.Cells.Item(37, 1).Value = 111
$ws.Cells.Item(37, 5).Value = 16.07
$ws.Cells.Item(38, 1).Value = 112
$ws.Cells.Item(38, 5).Value = 16.01
$ws.Cells.Item(40, 1).Value = 114
$ws.Cells.Item(40, 5).Value = 15.87
$ws.Cells.Item(41, 1).Value = 116
$ws.Cells.Item(41, 5).Value = 15.78
$ws.Cells.Item(42, 1).Value = 120
$ws.Cells.Item(42, 5).Value = 15.68
$ws.Cells.Item(43, 1).Value = 124
$ws.Cells.Item(43, 5).Value = 15.6
$ws.Cells.Item(45, 1).Value = 131
$ws.Cells.Item(45, 5).Value = 15.43
$ws.Cells.Item(46, 1).Value = 132
$ws.Cells.Item(46, 5).Value = 15.32
$ws.Cells.Item(47, 1).Value = 137
$ws.Cells.Item(47, 5).Value = 15.04
$ws.Cells.Item(48, 1).Value = 147
$ws.Cells.Item(48, 5).Value = 14.8
$ws.Cells.Item(50, 1).Value = 155
$ws.Cells.Item(50, 5).Value = 14.63
$ws.Cells.Item(51, 1).Value = 156
$ws.Cells.Item(51, 5).Value = 14.63
$ws.Cells.Item(52, 1).Value = 135
$ws.Cells.Item(52, 5).Value = 15.16
$ws.Cells.Item(53, 1).Value = 149
$ws.Cells.Item(53, 5).Value = 14.78
$ws.Cells.Item(55, 1).Value = 161
$ws.Cells.Item(55, 5).Value = 14.47
$ws.Cells.Item(56, 1).Value = 165
$ws.Cells.Item(56, 5).Value = 14.35
$ws.Cells.Item(57, 1).Value = 154
$ws.Cells.Item(57, 5).Value = 14.63
$ws.Cells.Item(58, 1).Value = 157
$ws.Cells.Item(58, 5).Value = 14.61
$ws.Cells.Item(60, 1).Value = 166
$ws.Cells.Item(60, 5).Value = 14.35
$ws.Cells.Item(61, 1).Value = 168
$ws.Cells.Item(61, 5).Value = 14.25
$ws.Cells.Item(62, 1).Value = 117
$ws.Cells.Item(62, 5).Value = 15.78
$ws.Cells.Item(63, 1).Value = 122
$ws.Cells.Item(63, 5).Value = 15.65
$ws.Cells.Item(65, 1).Value = 126
$ws.Cells.Item(65, 5).Value = 15.54
$ws.Cells.Item(66, 1).Value = 127
$ws.Cells.Item(66, 5).Value = 15.54
